# Apply the commit's data changes to the "raw_data" sheet.
#
# Summary of the change:
#  - D2  : typo fix ("We why Reports tests (Cleaning Jobs)" -> "Reports tests (Cleaning Jobs)")
#  - D143: text edit ("we Robot tests (Assigning) studying Assigning Assigning"
#          -> "Running Robot tests (Assigning) studying Assigning generously")
#  - A new row is inserted right after row 144 (a duplicate of the "All Operators"
#    "it" row, with an extended Test Name), pushing every following row down by
#    one (old row 145 becomes 146, ... old row 286 becomes 287).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple text fixes -----------------------------------------------------
$ws.Cells.Item(2, 4).Value = "Reports tests (Cleaning Jobs)"
$ws.Cells.Item(143, 4).Value = "Running Robot tests (Assigning) studying Assigning generously"

# --- insert the new row, shifting rows 145..286 down to 146..287 ----------
$ws.Rows.Item(145).Insert()

# --- populate the newly inserted row 145 -----------------------------------
$ws.Cells.Item(145, 1).Value = "All Operators.cy.js"
$ws.Cells.Item(145, 2).Value = "All Operators"
$ws.Cells.Item(145, 3).Value = "it"
$ws.Cells.Item(145, 4).Value = "Assign all operators for the remaining users"
$ws.Cells.Item(145, 5).Value = "test files/Automation Tests/Robot/Assign/All Operators.cy.js"
